# Trading update: 2026-02-17 15:24:24
# Append the newest MarketMaking trade (Trade #40, still OPEN) as row 41
# on both the "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A41").Value = 40

    # "2026-02-17" looks like a date, so a plain assignment gets auto-converted
    # to a date serial by the smart-parser. Force literal text via a leading
    # apostrophe, then strip the resulting quote-prefix style back to Normal
    # so the cell keeps the workbook's default (unstyled) formatting.
    $ws.Range("B41").Value = "'2026-02-17"
    $ws.Range("B41").Style = "Normal"

    $ws.Range("C41").Value = "15:23:46"
    $ws.Range("D41").Value = "MarketMaking"
    $ws.Range("E41").Value = "DOWN"
    $ws.Range("F41").Value = 0.48

    # Exit Price is blank (trade still OPEN) but stored as an empty string
    # cell rather than a truly empty cell, matching the source data export.
    $ws.Range("G41").Value = "'"
    $ws.Range("G41").Style = "Normal"

    $ws.Range("H41").Value = "OPEN"
    $ws.Range("I41").Value = 0
    $ws.Range("J41").Value = 0
    $ws.Range("K41").Value = 99.73988006373717
    $ws.Range("L41").Value = 0
    $ws.Range("M41").Value = 0
    $ws.Range("N41").Value = 0.6
    $ws.Range("O41").Value = "Normal spread capture: 19600 bps"

    # Exit Reason is also blank (trade still OPEN).
    $ws.Range("P41").Value = "'"
    $ws.Range("P41").Style = "Normal"

    $ws.Range("Q41").Value = 0
}
